{"js": "// This script replaces the date line and every \"a\u00f7b=c, d\" division\n// answer in the document's single table with its updated value, per the\n// commit diff. Every old text value below is unique within the document,\n// so a body-wide search-and-replace by exact text is unambiguous.\nconst replacements = [\n  [\"2025-07-16 Wednesday\", \"2025-07-17 Thursday\"],\n  [\"235\u00f79=26, 1\", \"871\u00f73=290, 1\"],\n  [\"780\u00f76=130, 0\", \"773\u00f76=128, 5\"],\n  [\"469\u00f79=52, 1\", \"693\u00f73=231, 0\"],\n  [\"828\u00f73=276, 0\", \"213\u00f72=106, 1\"],\n  [\"924\u00f78=115, 4\", \"358\u00f78=44, 6\"],\n  [\"539\u00f78=67, 3\", \"973\u00f79=108, 1\"],\n  [\"823\u00f76=137, 1\", \"293\u00f78=36, 5\"],\n  [\"896\u00f79=99, 5\", \"828\u00f72=414, 0\"],\n  [\"105\u00f75=21, 0\", \"386\u00f75=77, 1\"],\n  [\"214\u00f76=35, 4\", \"403\u00f75=80, 3\"],\n  [\"312\u00f78=39, 0\", \"132\u00f76=22, 0\"],\n  [\"781\u00f72=390, 1\", \"203\u00f79=22, 5\"],\n  [\"879\u00f78=109, 7\", \"986\u00f77=140, 6\"],\n  [\"123\u00f73=41, 0\", \"445\u00f76=74, 1\"],\n  [\"291\u00f76=48, 3\", \"841\u00f73=280, 1\"],\n  [\"120\u00f76=20, 0\", \"935\u00f73=311, 2\"],\n  [\"153\u00f73=51, 0\", \"235\u00f72=117, 1\"],\n  [\"463\u00f73=154, 1\", \"208\u00f75=41, 3\"],\n  [\"670\u00f75=134, 0\", \"689\u00f76=114, 5\"],\n  [\"341\u00f77=48, 5\", \"598\u00f77=85, 3\"],\n  [\"899\u00f75=179, 4\", \"899\u00f72=449, 1\"],\n  [\"887\u00f72=443, 1\", \"604\u00f74=151, 0\"],\n  [\"796\u00f75=159, 1\", \"711\u00f78=88, 7\"],\n  [\"846\u00f74=211, 2\", \"154\u00f77=22, 0\"],\n  [\"299\u00f76=49, 5\", \"424\u00f74=106, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every \"a\u00f7b=c, d\" division answer in the\n# document's single table with its updated value, per the commit diff.\n# Every \"old\" value below is unique within the document, so a\n# document-wide Find/Replace (wdReplaceAll) for each exact string is\n# unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2025-07-16 Wednesday', '2025-07-17 Thursday'),\n    @('235\u00f79=26, 1', '871\u00f73=290, 1'),\n    @('780\u00f76=130, 0', '773\u00f76=128, 5'),\n    @('469\u00f79=52, 1', '693\u00f73=231, 0'),\n    @('828\u00f73=276, 0', '213\u00f72=106, 1'),\n    @('924\u00f78=115, 4', '358\u00f78=44, 6'),\n    @('539\u00f78=67, 3', '973\u00f79=108, 1'),\n    @('823\u00f76=137, 1', '293\u00f78=36, 5'),\n    @('896\u00f79=99, 5', '828\u00f72=414, 0'),\n    @('105\u00f75=21, 0', '386\u00f75=77, 1'),\n    @('214\u00f76=35, 4', '403\u00f75=80, 3'),\n    @('312\u00f78=39, 0', '132\u00f76=22, 0'),\n    @('781\u00f72=390, 1', '203\u00f79=22, 5'),\n    @('879\u00f78=109, 7', '986\u00f77=140, 6'),\n    @('123\u00f73=41, 0', '445\u00f76=74, 1'),\n    @('291\u00f76=48, 3', '841\u00f73=280, 1'),\n    @('120\u00f76=20, 0', '935\u00f73=311, 2'),\n    @('153\u00f73=51, 0', '235\u00f72=117, 1'),\n    @('463\u00f73=154, 1', '208\u00f75=41, 3'),\n    @('670\u00f75=134, 0', '689\u00f76=114, 5'),\n    @('341\u00f77=48, 5', '598\u00f77=85, 3'),\n    @('899\u00f75=179, 4', '899\u00f72=449, 1'),\n    @('887\u00f72=443, 1', '604\u00f74=151, 0'),\n    @('796\u00f75=159, 1', '711\u00f78=88, 7'),\n    @('846\u00f74=211, 2', '154\u00f77=22, 0'),\n    @('299\u00f76=49, 5', '424\u00f74=106, 0'),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
